$wb = $excel.ActiveWorkbook

# Add a new worksheet "DocumentTypes" after the last existing sheet (SchoolDepartment)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "DocumentTypes"

# Fill the data column by column (A, then B, then C, then D)
$ws.Range("A1").Value = "ID card"
$ws.Range("A2").Value = "Passport"
$ws.Range("A3").Value = "Driver's license"

$ws.Range("B1").Value = "is required"
$ws.Range("B2").Value = "is not required"
$ws.Range("B3").Value = "is required"

$ws.Range("C1").Value = "employment contract"
$ws.Range("C2").Value = "rental contract"
$ws.Range("C3").Value = "vehicle registration"

$ws.Range("D1").Value = "is not required"
$ws.Range("D2").Value = "is  required"
$ws.Range("D3").Value = "is not required"

# Column widths matching the target layout
$ws.Columns.Item(1).ColumnWidth = 30.5
$ws.Columns.Item(2).ColumnWidth = 20.3
$ws.Columns.Item(3).ColumnWidth = 24
$ws.Columns.Item(4).ColumnWidth = 13.1666666666667

# Select D2 as the active cell on the new sheet, making it the selected/visible tab
$ws.Range("D2").Select()
